# Add the 5 May 2020 (serial 43956) data row to the "data" sheet's Table3.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$tbl = $ws.ListObjects.Item("Table3")

# Grow the table by one row; Excel places it right after the current last row.
$newListRow = $tbl.ListRows.Add()

$headerRow = $tbl.Range.Row
$newRow = $headerRow + $tbl.ListRows.Count

$ws.Cells.Item($newRow, 1).Value = 43956
$ws.Cells.Item($newRow, 2).Value = 33283
$ws.Cells.Item($newRow, 3).Value = 1832
$ws.Cells.Item($newRow, 4).Value = 59
$ws.Cells.Item($newRow, 5).Value = 5119

$ws.Cells.Item($newRow, 5).Select() | Out-Null
